# Updates the Step3_DataPts_* sheets to reflect the new
# "zero_before_threshold" behaviour: dims occurring before the
# noise_threshold / First Rise Point are now treated as 0, which shifts
# First_Noticeable_Increase_Index (col C), its cumulative value (col E)
# and the resulting Pulse_Width (col G) for each signal segment row.

$wb = $excel.ActiveWorkbook

$updates = @{
    "Step3_DataPts_0.5" = @{
        2 = @{ C = 87; E = 0.01428337551834045;  G = 44 }
        3 = @{ C = 48; E = 0.005494772194819301; G = 50 }
        4 = @{ C = 48; E = 0.005674122027836288; G = 50 }
        5 = @{ C = 87; E = 0.008391815036828829; G = 45 }
        6 = @{ C = 47; E = 0.004997065866706255; G = 50 }
    }
    "Step3_DataPts_0.7" = @{
        2 = @{ C = 87; E = 0.01428337551834045;  G = 54 }
        3 = @{ C = 48; E = 0.005494772194819301; G = 61 }
        4 = @{ C = 48; E = 0.005674122027836288; G = 61 }
        5 = @{ C = 87; E = 0.008391815036828829; G = 54 }
        6 = @{ C = 47; E = 0.004997065866706255; G = 61 }
    }
    "Step3_DataPts_0.8" = @{
        2 = @{ C = 87; E = 0.01428337551834045;  G = 61 }
        3 = @{ C = 48; E = 0.005494772194819301; G = 78 }
        4 = @{ C = 48; E = 0.005674122027836288; G = 77 }
        5 = @{ C = 87; E = 0.008391815036828829; G = 62 }
        6 = @{ C = 47; E = 0.004997065866706255; G = 77 }
    }
    "Step3_DataPts_0.9" = @{
        2 = @{ C = 87; E = 0.01428337551834045;  G = 80 }
        3 = @{ C = 48; E = 0.005494772194819301; G = 112 }
        4 = @{ C = 48; E = 0.005674122027836288; G = 111 }
        5 = @{ C = 87; E = 0.008391815036828829; G = 80 }
        6 = @{ C = 47; E = 0.004997065866706255; G = 112 }
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $cols = $rows[$rowNum]
        foreach ($colLetter in $cols.Keys) {
            $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
        }
    }
}
